# Insert a new weekly price-record row at position 146 (pushing the
# existing rows 146-211 down to 147-212), matching the author's commit
# "Fruta / hortaliza, semanal" which adds one more weekly data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 146, shifting rows down.
$ws.Rows.Item(146).Insert()

# Populate the newly inserted row 146 with the new record's data.
$ws.Cells.Item(146, 1).Value  = 2
$ws.Cells.Item(146, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(146, 3).Value  = "Coquimbo"
$ws.Cells.Item(146, 4).Value  = 45029
$ws.Cells.Item(146, 5).Value  = 4
$ws.Cells.Item(146, 6).Value  = 100112043
$ws.Cells.Item(146, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(146, 8).Value  = "Sin especificar"
$ws.Cells.Item(146, 9).Value  = "Primera"
$ws.Cells.Item(146, 10).Value = 500
$ws.Cells.Item(146, 11).Value = 7000
$ws.Cells.Item(146, 12).Value = 7500
$ws.Cells.Item(146, 13).Value = 7250
$ws.Cells.Item(146, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(146, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(146, 16).Value = 104
$ws.Cells.Item(146, 17).Value = 70
$ws.Cells.Item(146, 18).Value = "Hortaliza"
